# Add two new rows to the end of the data table (Sheet1), duplicating the
# values of the last existing row but advancing the date (column A) by one
# day for each new row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Locate the last populated row in column A (xlUp = -4162).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$lastCol = 10   # columns A..J

$newRowsCount = 2

for ($i = 1; $i -le $newRowsCount; $i++) {
    $srcRow = $lastRow + $i - 1
    $dstRow = $lastRow + $i

    $src = $ws.Range($ws.Cells.Item($srcRow, 1), $ws.Cells.Item($srcRow, $lastCol))
    $dst = $ws.Range($ws.Cells.Item($dstRow, 1), $ws.Cells.Item($dstRow, $lastCol))

    # Copy values + formatting (style) from the row above into the new row.
    $src.Copy($dst)

    # Advance the date serial in column A by one day relative to the row
    # that was just copied from. Value2 returns the raw numeric (serial
    # date) representation rather than a formatted DateTime, so arithmetic
    # keeps the cell numeric instead of turning it into text.
    $prevDateSerial = $ws.Cells.Item($srcRow, 1).Value2()
    $ws.Cells.Item($dstRow, 1).Value = $prevDateSerial + 1
}
